# NATMI re-run with new TPM-derived ligand/receptor expression values.
# The FAPs/MuSCs sending-cluster block replaces the old ECs block, and the
# previous MuSCs->{ECs,FAPs,MuSCs} rows (old rows 8-10) are dropped entirely,
# shrinking the table from 9 data rows (A2:T10) down to 6 data rows (A2:T7).

function Set-FullRow($ws, $row, $values) {
    # Build a proper 1 x N 2-D object array (required for a single Range.Value2
    # assignment to fan out across multiple columns in one shot).
    $arr = New-Object 'object[,]' 1,$values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $rng = "A" + $row + ":T" + $row
    $ws.Range($rng).Value2 = $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expr value,
# Ligand total expr value, Ligand derived specificity (avg), Ligand derived
# specificity (total), Receptor-expressing cells, Receptor detection rate,
# Receptor average expr value, Receptor total expr value, Receptor derived
# specificity (avg), Receptor derived specificity (total), Edge average expr
# weight, Edge total expr weight, Edge average expr derived specificity,
# Edge total expr derived specificity.

Set-FullRow $ws 2 @("FAPs","Angpt1","Itga5","ECs",3,1,11.05178533333333,33.155356,0.9017494976312432,0.9017494976312432,3,1,10.34761366666667,31.042841,0.2299953477621856,0.2299953477621856,114.3596049562662,1029.236444606396,0.2073981893020739,0.207398189302074)
Set-FullRow $ws 3 @("FAPs","Angpt1","Itga5","FAPs",3,1,11.05178533333333,33.155356,0.9017494976312432,0.9017494976312432,3,1,30.56986233333333,91.709587,0.6794731949692173,0.6794731949692174,337.8515561775524,3040.664005597972,0.6127146122173875,0.6127146122173877)
Set-FullRow $ws 4 @("FAPs","Angpt1","Itga5","MuSCs",3,1,11.05178533333333,33.155356,0.9017494976312432,0.9017494976312432,3,1,4.073058666666666,12.219176,0.09053145726859702,0.09053145726859703,45.01457003407288,405.131130306656,0.08163669611178172,0.08163669611178173)
Set-FullRow $ws 5 @("MuSCs","Angpt1","Itga5","ECs",3,1,1.204152,3.612456,0.09825050236875665,0.09825050236875667,3,1,10.34761366666667,31.042841,0.2299953477621856,0.2299953477621856,12.460099691944,112.140897227496,0.02259715846011162,0.02259715846011163)
Set-FullRow $ws 6 @("MuSCs","Angpt1","Itga5","FAPs",3,1,1.204152,3.612456,0.09825050236875665,0.09825050236875667,3,1,30.56986233333333,91.709587,0.6794731949692173,0.6794731949692174,36.810760868408,331.296847815672,0.06675858275182973,0.06675858275182976)
Set-FullRow $ws 7 @("MuSCs","Angpt1","Itga5","MuSCs",3,1,1.204152,3.612456,0.09825050236875665,0.09825050236875667,3,1,4.073058666666666,12.219176,0.09053145726859702,0.09053145726859703,4.904581739583999,44.14123565625599,0.008894761156815282,0.008894761156815286)

# Drop the trailing rows that no longer exist in the refreshed TPM output
# (old rows 8, 9 and 10), shrinking the sheet dimension to A1:T7.
$ws.Range("A8:A10").EntireRow.Delete() | Out-Null

Write-Host "Updated sheet to new TPM values; used range now" $ws.UsedRange.Address()
